$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Feb 24 23:03:22 EST 2025"
$ws.Range("B3").Value = "Mon Feb 24 23:03:37 EST 2025"
$ws.Range("B4").Value = "Mon Feb 24 23:03:51 EST 2025"
$ws.Range("B5").Value = "Mon Feb 24 23:04:06 EST 2025"
$ws.Range("B6").Value = "Mon Feb 24 23:04:21 EST 2025"
$ws.Range("B7").Value = "Mon Feb 24 23:04:35 EST 2025"
